$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "SFH"

$rows = @(
  @(2, 1600, 40802.01931231975),
  @(3, 1601, 49.15613759094706),
  @(4, 1602, 49.9545662217164),
  @(5, 1603, 50.77057599935205),
  @(6, 1604, 51.60477224137885),
  @(7, 1605, 52.45778484216541),
  @(8, 1606, 53.33026900979727),
  @(9, 1607, 54.2229059976972),
  @(10, 1608, 55.13640382852392),
  @(11, 1609, 56.0714980070085),
  @(12, 1610, 57.02895221948558),
  @(13, 1611, 58.00955901712577),
  @(14, 1612, 59.01414047896407),
  @(15, 1613, 60.04354885243271),
  @(16, 1614, 61.0986671668802),
  @(17, 1615, 62.18040981739722),
  @(18, 1616, 63.28972311433482),
  @(19, 1617, 64.42758579514161),
  @(20, 1618, 65.5950094943567),
  @(21, 1619, 66.79303916718354),
  @(22, 1620, 68.02275346310235),
  @(23, 1621, 69.28526504426854),
  @(24, 1622, 70.58172084471097),
  @(25, 1623, 71.91330226575559),
  @(26, 1624, 73.28122530246944),
  @(27, 1625, 74.68674059711546),
  @(28, 1626, 76.13113341421196),
  @(29, 1627, 77.61572353267792),
  @(30, 1628, 79.14186505038492),
  @(31, 1629, 80.71094609581388),
  @(32, 1630, 82.32438844251651),
  @(33, 1631, 83.98364702142709),
  @(34, 1632, 85.69020932627204),
  @(35, 1633, 87.44559470773657),
  @(36, 1634, 89.25135355163569),
  @(37, 1635, 91.10906633713599),
  @(38, 1636, 93.02034257035662),
  @(39, 1637, 94.98681959004901),
  @(40, 1638, 97.01016124113579),
  @(41, 1639, 99.09205641285683),
  @(42, 1640, 101.234217438021),
  @(43, 1641, 103.4383783512978),
  @(44, 1642, 105.7062930027724),
  @(45, 1643, 108.0397330257971),
  @(46, 1644, 110.4404856565631),
  @(47, 1645, 112.9103514041932),
  @(48, 1646, 115.4511415710698),
  @(49, 1647, 118.0646756219653),
  @(50, 1648, 120.7527784028512),
  @(51, 1649, 123.5172772099776),
  @(52, 1650, 126.3599987099635),
  @(53, 1651, 129.2827657132635),
  @(54, 1652, 132.2873938028891),
  @(55, 1653, 135.3756878222032),
  @(56, 1654, 138.5494382244977),
  @(57, 1655, 141.8104172895173),
  @(58, 1656, 145.1603752120231),
  @(59, 1657, 148.6010360673725),
  @(60, 1658, 152.1340936615071),
  @(61, 1659, 155.7612072722547),
  @(62, 1660, 159.4839972897009),
  @(63, 1661, 163.304040764367),
  @(64, 1662, 167.2228668732244),
  @(65, 1663, 171.2419523123242),
  @(66, 1664, 175.36271662871),
  @(67, 1665, 179.5865175013685),
  @(68, 1666, 183.9146459847571),
  @(69, 1667, 188.3483217274562),
  @(70, 1668, 192.8886881792907),
  @(71, 1669, 197.5368078012591),
  @(72, 1670, 202.2936572935454),
  @(73, 1671, 207.1601228560786),
  @(74, 1672, 212.1369954984329),
  @(75, 1673, 217.2249664148821),
  @(76, 1674, 222.4246224416177),
  @(77, 1675, 227.7364416132594),
  @(78, 1676, 233.1607888363863),
  @(79, 1677, 238.6979116974678),
  @(80, 1678, 244.3479364237282),
  @(81, 1679, 250.1108640146885),
  @(82, 1680, 255.9865665628374),
  @(83, 1681, 261.974783781636),
  @(84, 1682, 268.0751197591221),
  @(85, 1683, 274.2870399550441),
  @(86, 1684, 280.6098684596459),
  @(87, 1685, 287.0427855309827),
  @(88, 1686, 293.5848254290194),
  @(89, 1687, 300.2348745625997),
  @(90, 1688, 306.9916699649411),
  @(91, 1689, 313.8537981146534),
  @(92, 1690, 320.8196941160328),
  @(93, 1691, 327.8876412528793),
  @(94, 1692, 335.0557709295949),
  @(95, 1693, 342.3220630115276),
  @(96, 1694, 349.6843465755453),
  @(97, 1695, 357.1403010816719),
  @(98, 1696, 364.6874579742936),
  @(99, 1697, 372.32320272097),
  @(100, 1698, 380.0447772952478),
  @(101, 1699, 387.8492831082889),
  @(102, 1700, 395.7336843938841),
  @(103, 1701, 403.6948120479448),
  @(104, 1702, 411.7293679236161),
  @(105, 1703, 419.8339295818358),
  @(106, 1704, 428.0049554933034),
  @(107, 1705, 436.2387906894068),
  @(108, 1706, 444.531672854869),
  @(109, 1707, 452.8797388557571),
  @(110, 1708, 461.2790316931618),
  @(111, 1709, 469.7255078716404),
  @(112, 1710, 478.2150451704071),
  @(113, 1711, 486.7434508023447),
  @(114, 1712, 495.306469945878),
  @(115, 1713, 503.8997946315054),
  @(116, 1714, 512.519072964319),
  @(117, 1715, 521.1599186618876),
  @(118, 1716, 529.8179208855581),
  @(119, 1717, 538.4886543412333),
  @(120, 1718, 547.1676896253562),
  @(121, 1719, 555.8506037895396),
  @(122, 1720, 564.5329910967205),
  @(123, 1721, 573.2104739401559),
  @(124, 1722, 581.8787138959538),
  @(125, 1723, 590.5334228786169),
  @(126, 1724, 599.1703743684682),
  @(127, 1725, 607.7854146787248),
  @(128, 1726, 616.3744742302189),
  @(129, 1727, 624.9335787999913),
  @(130, 1728, 633.4588607114741),
  @(131, 1729, 641.9465699312346),
  @(132, 1730, 650.3930850402589),
  @(133, 1731, 658.79492404506),
  @(134, 1732, 667.1487549950963),
  @(135, 1733, 675.451406374778),
  @(136, 1734, 683.699877235106),
  @(137, 1735, 691.8913470352505),
  @(138, 1736, 700.023185161156),
  @(139, 1737, 708.0929600908128),
  @(140, 1738, 716.0984481775139),
  @(141, 1739, 724.0376420216618),
  @(142, 1740, 731.9087584041969),
  @(143, 1741, 739.7102457563651),
  @(144, 1742, 747.4407911401847),
  @(145, 1743, 755.099326718073),
  @(146, 1744, 762.6850356890137),
  @(147, 1745, 770.1973576733842),
  @(148, 1746, 777.6359935269502),
  @(149, 1747, 785.0009095699531),
  @(150, 1748, 792.2923412167936),
  @(151, 1749, 799.5107959940582),
  @(152, 1750, 806.6570559387515),
  @(153, 1751, 813.7321793675812),
  @(154, 1752, 820.7375020135836),
  @(155, 1753, 827.674637525709),
  @(156, 1754, 834.5454773310838),
  @(157, 1755, 841.3521898619807),
  @(158, 1756, 848.0972191492576),
  @(159, 1757, 854.7832827903218),
  @(160, 1758, 861.4133692976908),
  @(161, 1759, 867.9907348394729),
  @(162, 1760, 874.5188993834787),
  @(163, 1761, 881.0016422594097),
  @(164, 1762, 887.442997155699),
  @(165, 1763, 893.8472465683911),
  @(166, 1764, 900.218915723335),
  @(167, 1765, 906.5627659929604),
  @(168, 1766, 912.8837878308974),
  @(169, 1767, 919.187193250106),
  @(170, 1768, 925.4784078712037),
  @(171, 1769, 931.7630625684329),
  @(172, 1770, 938.046984743048),
  @(173, 1771, 944.3361892543309),
  @(174, 1772, 950.6368690397472),
  @(175, 1773, 956.9553854562389),
  @(176, 1774, 963.2982583761785),
  @(177, 1775, 969.6721560714757),
  @(178, 1776, 976.0838849198307),
  @(179, 1777, 982.5403789680879),
  @(180, 1778, 989.0486893870142),
  @(181, 1779, 995.615973852853),
  @(182, 1780, 1002.249485890319),
  @(183, 1781, 1008.956564211514),
  @(184, 1782, 1015.744622085251),
  @(185, 1783, 1022.621136771328),
  @(186, 1784, 1029.593639052244),
  @(187, 1785, 1036.669702895728),
  @(188, 1786, 1043.856935280467),
  @(189, 1787, 1051.162966214888),
  @(190, 1788, 1058.595438981091),
  @(191, 1789, 1066.162000631288),
  @(192, 1790, 1073.870292765398),
  @(193, 1791, 1081.72794261756),
  @(194, 1792, 1089.74255447483),
  @(195, 1793, 1097.921701453835),
  @(196, 1794, 1106.272917657652),
  @(197, 1795, 1114.803690733424),
  @(198, 1796, 1123.521454850677),
  @(199, 1797, 1132.433584118545),
  @(200, 1798, 1141.547386457454),
  @(201, 1799, 1150.870097941007),
  @(202, 1800, 1160.408877620377),
  @(203, 1801, 1170.171063960404),
  @(204, 1802, 1180.163155112692),
  @(205, 1803, 1190.392287424238),
  @(206, 1804, 1200.865270253245),
  @(207, 1805, 1211.588817644456),
  @(208, 1806, 1222.569546230439),
  @(209, 1807, 1233.813973660908),
  @(210, 1808, 1245.328517559545),
  @(211, 1809, 1257.119495007001),
  @(212, 1810, 1269.193122547158),
  @(213, 1811, 1281.55551671175),
  @(214, 1812, 1294.212695058007),
  @(215, 1813, 1307.170577712723),
  @(216, 1814, 1320.434989413106),
  @(217, 1815, 1334.011662035959),
  @(218, 1816, 1347.906237604498),
  @(219, 1817, 1362.124271760324),
  @(220, 1818, 1376.671237688548),
  @(221, 1819, 1391.552530483511),
  @(222, 1820, 1406.773471936802),
  @(223, 1821, 1422.339315738072),
  @(224, 1822, 1438.255253067018),
  @(225, 1823, 1454.526418565046),
  @(226, 1824, 1471.157896663636),
  @(227, 1825, 1488.154728258994),
  @(228, 1826, 1505.521917707465),
  @(229, 1827, 1523.264440128415),
  @(230, 1828, 1541.38724899303),
  @(231, 1829, 1559.895283981973),
  @(232, 1830, 1578.793479089268),
  @(233, 1831, 1598.086770956322),
  @(234, 1832, 1617.780107415122),
  @(235, 1833, 1637.878456219146),
  @(236, 1834, 1658.386813945238),
  @(237, 1835, 1679.310215044075),
  @(238, 1836, 1700.653741020251),
  @(239, 1837, 1722.422529722709),
  @(240, 1838, 1744.621784725047),
  @(241, 1839, 1767.256784775204),
  @(242, 1840, 1790.332893297522),
  @(243, 1841, 1813.855567924855),
  @(244, 1842, 1837.830370041942),
  @(245, 1843, 1862.262974322586),
  @(246, 1844, 1887.159178238036),
  @(247, 1845, 1912.524911520732),
  @(248, 1846, 1938.366245561557),
  @(249, 1847, 1964.689402722904),
  @(250, 1848, 1991.500765547927),
  @(251, 1849, 2018.806885847525),
  @(252, 1850, 2046.614493643931),
  @(253, 1851, 2074.930505954591),
  @(254, 1852, 2103.762035394938),
  @(255, 1853, 2133.116398579951),
  @(256, 1854, 2163.001124307367),
  @(257, 1855, 2193.423961500781),
  @(258, 1856, 2224.392886892285),
  @(259, 1857, 2255.916112424934),
  @(260, 1858, 2288.002092355164),
  @(261, 1859, 2320.659530030531),
  @(262, 1860, 2353.897384326202),
  @(263, 1861, 2387.724875712958),
  @(264, 1862, 2422.151491936378),
  @(265, 1863, 2457.186993287084),
  @(266, 1864, 2492.841417429135),
  @(267, 1865, 2529.125083773593),
  @(268, 1866, 2566.048597363488),
  @(269, 1867, 2603.622852246337),
  @(270, 1868, 2641.859034310387),
  @(271, 1869, 2680.768623552637),
  @(272, 1870, 2720.36339575455),
  @(273, 1871, 2760.6554235323),
  @(274, 1872, 2801.657076734699),
  @(275, 1873, 2843.381022156539),
  @(276, 1874, 2885.840222534595),
  @(277, 1875, 2929.047934796541),
  @(278, 1876, 2973.017707526045),
  @(279, 1877, 3017.763377612099),
  @(280, 1878, 3063.299066045798),
  @(281, 1879, 3109.639172830436),
  @(282, 1880, 3156.798370966535),
  @(283, 1881, 3204.791599475737),
  @(284, 1882, 3253.634055424926),
  @(285, 1883, 3303.3411849125),
  @(286, 1884, 3353.928672978856),
  @(287, 1885, 3405.412432398537),
  @(288, 1886, 3457.808591319444),
  @(289, 1887, 3511.133479705347),
  @(290, 1888, 3565.403614544694),
  @(291, 1889, 3620.635683785341),
  @(292, 1890, 3676.846528956005),
  @(293, 1891, 3734.053126436269),
  @(294, 1892, 3792.272567339652),
  @(295, 1893, 3851.522035967085),
  @(296, 1894, 3911.818786802785),
  @(297, 1895, 3973.180120013421),
  @(298, 1896, 4035.623355419353),
  @(299, 1897, 4099.165804907661),
  @(300, 1898, 4163.824743261097),
  @(301, 1899, 4229.617377374315),
  @(302, 1900, 4296.560813834873),
  @(303, 1901, 4364.672024856191),
  @(304, 1902, 4433.967812538945),
  @(305, 1903, 4504.464771453842),
  @(306, 1904, 4576.179249541552),
  @(307, 1905, 4649.127307323627),
  @(308, 1906, 4723.324675429556),
  @(309, 1907, 4798.78671045367),
  @(310, 1908, 4875.528349153838),
  @(311, 1909, 4953.564061012948),
  @(312, 1910, 5032.907799201491),
  @(313, 1911, 5113.572949973648),
  @(314, 1912, 5195.572280541088),
  @(315, 1913, 5278.917885486871),
  @(316, 1914, 5363.62113177951),
  @(317, 1915, 5449.692602459963),
  @(318, 1916, 5537.142039093213),
  @(319, 1917, 5625.978283070373),
  @(320, 1918, 5716.209215872871),
  @(321, 1919, 5807.841698420432),
  @(322, 1920, 5900.881509620703),
  @(323, 1921, 5995.333284275597),
  @(324, 1922, 6091.200450484396),
  @(325, 1923, 6188.485166725873),
  @(326, 1924, 6287.188258782327),
  @(327, 1925, 6387.309156707593),
  @(328, 1926, 6488.845832041718),
  @(329, 1927, 6591.794735485224),
  @(330, 1928, 6696.150735266089),
  @(331, 1929, 6801.907056440094),
  @(332, 1930, 6909.055221378264),
  @(333, 1931, 7017.584991702349),
  @(334, 1932, 7127.48431195016),
  @(335, 1933, 7238.739255250446),
  @(336, 1934, 7351.333971305094),
  @(337, 1935, 7465.250636983444),
  @(338, 1936, 7580.469409836267),
  @(339, 1937, 7696.968384846134),
  @(340, 1938, 7814.723554746885),
  @(341, 1939, 7933.708774225192),
  @(342, 1940, 8053.895728341571),
  @(343, 1941, 8175.253905503288),
  @(344, 1942, 8297.750575312732),
  @(345, 1943, 8421.350771628751),
  @(346, 1944, 8546.017281150587),
  @(347, 1945, 8671.710637865279),
  @(348, 1946, 8798.389123643121),
  @(349, 1947, 8926.008775309976),
  @(350, 1948, 9054.523398459498),
  @(351, 1949, 9183.88458830691),
  @(352, 1950, 9314.041757824067),
  @(353, 1951, 9444.942173412921),
  @(354, 1952, 9576.530998328955),
  @(355, 1953, 9708.751344064229),
  @(356, 1954, 9841.544329858893),
  @(357, 1955, 9974.849150494243),
  @(358, 1956, 9648.422021059796),
  @(359, 1957, 9776.454292890461),
  @(360, 1958, 9904.790734282322),
  @(361, 1959, 10033.36676928167),
  @(362, 1960, 10162.11643019076),
  @(363, 1961, 10290.97247105825),
  @(364, 1962, 10419.86648923632),
  @(365, 1963, 10548.72905486466),
  @(366, 1964, 10677.48984813675),
  @(367, 1965, 10806.07780416475),
  @(368, 1966, 10934.42126518825),
  @(369, 1967, 11062.44813988589),
  @(370, 1968, 11190.08606945521),
  @(371, 1969, 11317.26260013278),
  @(372, 1970, 11443.90536176952),
  @(373, 1971, 12525.39438235054),
  @(374, 1972, 12661.10600145088),
  @(375, 1973, 12796.0072978233),
  @(376, 1974, 12930.02233462394),
  @(377, 1975, 13063.07631327711),
  @(378, 1976, 13195.09579534989),
  @(379, 1977, 13326.00892708885),
  @(380, 1978, 13455.745665947),
  @(381, 1979, 13584.2380084111),
  @(382, 1980, 13711.42021841431),
  @(383, 1981, 14792.39797354074),
  @(384, 1982, 14925.35838848857),
  @(385, 1983, 15056.7243914057),
  @(386, 1984, 15186.43783146891),
  @(387, 1985, 15314.44418470614),
  @(388, 1986, 15440.69278538259),
  @(389, 1987, 15565.13705245964),
  @(390, 1988, 15687.73471034337),
  @(391, 1989, 15808.44800313778),
  @(392, 1990, 15927.24390164987),
  @(393, 1991, 12477.9900996606),
  @(394, 1992, 12567.33733088284),
  @(395, 1993, 12655.13985136096),
  @(396, 1994, 12741.38815470402),
  @(397, 1995, 12826.07711111212),
  @(398, 1996, 12909.20608282247),
  @(399, 1997, 12990.77903017414),
  @(400, 1998, 13070.80460779751),
  @(401, 1999, 13149.29625050741),
  @(402, 2000, 13226.27224842314),
  @(403, 2001, 9407.887849263001),
  @(404, 2002, 9460.239234968991),
  @(405, 2003, 9511.578478675992),
  @(406, 2004, 9561.932621245875),
  @(407, 2005, 9611.332278873711),
  @(408, 2006, 9659.811644516838),
  @(409, 2007, 9707.408479760727),
  @(410, 2008, 9754.164096812417),
  @(411, 2009, 9800.123330391727),
  @(412, 2010, 9845.334499201914),
  @(413, 2011, 9991.874344193982),
  @(414, 2012, 10036.20062513969),
  @(415, 2013, 10079.93814241664),
  @(416, 2014, 10123.14918343805),
  @(417, 2015, 10165.89919311108),
  @(418, 2016, 10208.25666562707),
  @(419, 2017, 10250.2930235946),
  @(420, 2018, 10292.08248417142),
  @(421, 2019, 10333.70191185926),
  @(422, 2020, 10375.23065760291),
  @(423, 2021, 10416.75038387236),
  @(424, 2022, 10458.34487535952),
  @(425, 2023, 10500.09983499424),
  @(426, 2024, 10542.10266492985),
  @(427, 2025, 10584.44223223253),
  @(428, 2026, 10627.20861900341),
  @(429, 2027, 10670.49285669636),
  @(430, 2028, 10714.38664446927),
  @(431, 2029, 10758.98205143419),
  @(432, 2030, 10804.37120274763),
  @(433, 2031, 10850.64594957725),
  @(434, 2032, 10897.89752303819),
  @(435, 2033, 10946.21617233962),
  @(436, 2034, 10995.69078744448),
  @(437, 2035, 11046.40850673984),
  @(438, 2036, 11098.45431029006),
  @(439, 2037, 11151.91059947101),
  @(440, 2038, 11206.85676388627),
  @(441, 2039, 11263.36873671125),
  @(442, 2040, 11321.51853976431),
  @(443, 2041, 11381.37381981988),
  @(444, 2042, 11442.99737791391),
  @(445, 2043, 11506.44669359043),
  @(446, 2044, 11571.77344628323),
  @(447, 2045, 11639.02303625528),
  @(448, 2046, 11708.23410775656),
  @(449, 2047, 11779.43807731464),
  @(450, 2048, 11852.65867024238),
  @(451, 2049, 11927.9114687628),
  @(452, 2050, 12005.20347527854)
)

foreach ($r in $rows) {
  $ws.Cells.Item($r[0], 1).Value = $r[1]
  $ws.Cells.Item($r[0], 2).Value = $r[2]
}
